$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consent")
$ws.Range("B2:D2").WrapText = $true
